$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "Recorded By" email lists on rows 2 and 3 ---
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value = "hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg"

# --- Row 10 (HISTOLOGY, Year 2 / C1, session 2) moved from "Pending" to "Not Recorded" ---
# Copy the existing "Not Recorded" row formatting (row 29) onto row 10 so the
# same fill/font combination gets reused, then update the status text.
$ws.Range("A29:I29").Copy() | Out-Null
$ws.Range("A10:I10").PasteSpecial(-4122) | Out-Null
$ws.Range("I10").Value = "Not Recorded"

# --- Update the dependent summary statistics (Missing/Pending sessions) ---
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 20
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 20

$excel.CutCopyMode = 0
